{"js": "// Assignment 2 groups roster update:\n// Add a new group's members (\"Paschal Chidiutor Ibeh, Oluwadamilola Ogundipe,\n// Abdulrahman Hamid\") to the bullet item that currently only holds a lone\n// \".\" placeholder \u2014 the bullet right after the \"Chloe Quijano, Ibeh\n// Mary-Anne, Fahima Nawshin.\" group.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the \"Chloe Quijano...\" group paragraph, then the still-empty\n// placeholder bullet that immediately follows it (a lone \".\").\nlet placeholderIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Chloe\") !== -1 && text.indexOf(\"Quijano\") !== -1) {\n    const candidate = paragraphs.items[i + 1];\n    if (candidate) {\n      placeholderIndex = i + 1;\n    }\n    break;\n  }\n}\n\nif (placeholderIndex === -1) {\n  throw new Error(\"Could not locate the placeholder group paragraph.\");\n}\n\nconst placeholder = paragraphs.items[placeholderIndex];\nplaceholder.load(\"text\");\nawait context.sync();\n\nif (placeholder.text === \".\") {\n  placeholder.insertText(\n    \"Paschal Chidiutor Ibeh, Oluwadamilola Ogundipe, Abdulrahman Hamid\",\n    \"Start\"\n  );\n  await context.sync();\n}\n", "ps1": "# Assignment 2 groups roster update:\n# Add a new group's members (\"Paschal Chidiutor Ibeh, Oluwadamilola Ogundipe,\n# Abdulrahman Hamid\") to the bullet item that currently only holds a lone\n# \".\" placeholder -- the bullet right after the \"Chloe Quijano, Ibeh\n# Mary-Anne, Fahima Nawshin.\" group.\n\n$d = $word.ActiveDocument\n\n$placeholderIndex = -1\n$i = 1\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Chloe*\" -and $t -like \"*Quijano*\") {\n        $placeholderIndex = $i + 1\n        break\n    }\n    $i++\n}\n\nif ($placeholderIndex -ne -1) {\n    $placeholder = $d.Paragraphs.Item($placeholderIndex).Range\n    # Paragraph.Range.Text includes the trailing paragraph mark (chr 13),\n    # so trim it off before comparing against the literal placeholder \".\".\n    $placeholderText = $placeholder.Text.TrimEnd([char]13, [char]10)\n    if ($placeholderText -eq \".\") {\n        $placeholder.InsertBefore(\"Paschal Chidiutor Ibeh, Oluwadamilola Ogundipe, Abdulrahman Hamid\")\n    }\n}\n"}
